$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 192, pushing existing rows 192-272 down to 193-273
$ws.Rows.Item(192).Insert()

# Populate the newly inserted row 192 with the new record
$ws.Range("A192").Value = 5
$ws.Range("B192").Value = "Macroferia Regional de Talca"
$ws.Range("C192").Value = "Maule"
$ws.Range("D192").Value = 44455
$ws.Range("E192").Value = 7
$ws.Range("F192").Value = "Fruta"
$ws.Range("G192").Value = 100102
$ws.Range("H192").Value = "Cítricos"
$ws.Range("I192").Value = 100102005
$ws.Range("J192").Value = "Naranja"
$ws.Range("K192").Value = "Navel Late"
$ws.Range("L192").Value = "Primera"
$ws.Range("M192").Value = 300
$ws.Range("N192").Value = 6000
$ws.Range("O192").Value = 6000
$ws.Range("P192").Value = 6000
$ws.Range("Q192").Value = '$/bandeja 15 kilos granel'
$ws.Range("R192").Value = "Región de O'Higgins"
$ws.Range("S192").Value = 400
$ws.Range("T192").Value = 15
